{"js": "// Apply the \"Level 1 mission complete\" edits to the home utilities contacts\n// document: update the intro heading, refresh provider descriptions /\n// addresses / phone numbers / websites for Sewer, Water and\n// Garbage/Recycling, rename the \"Additional Information\" heading, and swap\n// the old Hazardous Waste / Composting bullets for four new contact lines.\n\nconst body = context.document.body;\n\n// Simple text substitutions: each pair is matched literally (no wildcards)\n// and replaced in place. Search text is specific enough to be unique.\nconst replacements = [\n  [\n    \"<b>Utility Providers for San Jose, CA 95148</b>\",\n    \"<b>Utility Providers in San Jose, CA 95148</b>\",\n  ],\n  [\n    \"- Description: PG&E provides natural gas and electric service to approximately 16 million Californians.\",\n    \"- Description: PG&E provides natural gas and electric service to approximately 16 million people throughout a 70,000-square-mile service area in Northern and Central California.\",\n  ],\n  [\n    \"- Description: San Jose Water provides water service to over one million people in the greater San Jose metropolitan area.\",\n    \"- Description: San Jose Water provides essential water services to approximately 1 million people in the greater San Jose metropolitan area.\",\n  ],\n  [\n    \"- Address: 2130 The Alameda, San Jose, CA 95126\",\n    \"- Address: 110 Paseo De San Antonio, San Jose, CA 95112\",\n  ],\n  [\n    \"- Description: The City of San Jose provides wastewater services to residents and businesses within the city limits.\",\n    \"- Description: The City of San Jose's Environmental Services department is responsible for wastewater management, recycling, and solid waste services.\",\n  ],\n  [\n    \"- Address: 505 Los Coches St, San Jose, CA 95122\",\n    \"- Address: 5055 Almaden Expy, San Jose, CA 95118\",\n  ],\n  // The Sewer contact number/website line \u2014 matched together with the\n  // preceding address text so this occurrence (and not the identical phone\n  // number that also appears later, inside the Hazardous Waste bullet) is\n  // the one that gets updated.\n  [\n    \"- Contact Number: (408) 535-6000\\u000b- Website: <https://www.sanjoseca.gov/services/wastewater>\",\n    \"- Contact Number: (408) 277-4343\\u000b- Website: <https://www.sanjoseca.gov/index.aspx?NID=242>\",\n  ],\n  [\n    \"- Description: Republic Services provides residential and commercial waste and recycling services in San Jose.\",\n    \"- Description: Republic Services provides recycling, trash, and yard waste collection services to residents and businesses in San Jose.\",\n  ],\n  [\n    \"- Address: 1531 Oakland Rd, San Jose, CA 95110\",\n    \"- Address: 2550 S 10th St, San Jose, CA 95112\",\n  ],\n  [\n    \"- Contact Number: (408) 629-5000\",\n    \"- Contact Number: (408) 629-8500\",\n  ],\n  [\n    \"<b>Additional Information:</b>\",\n    \"<b>Additional Contact Information for Garbage/Recycling:</b>\",\n  ],\n  // Replace the two old bullet lines (Hazardous Waste + Composting,\n  // joined by the line break between them) with the four new contact lines.\n  [\n    \"- <b>Hazardous Waste:</b> The City of San Jose's Household Hazardous Waste Program provides safe disposal of hazardous materials. Contact (408) 535-6000 for more information.\\u000b- <b>Composting:</b> San Jose's GreenCycle Composting Program offers curbside composting services. Visit <https://www.sanjoseca.gov/services/greencycle> for details.\",\n    \"- Collection Schedule: <https://www.sanjoseca.gov/index.aspx?NID=242>\\u000b- Bulky Item Pickup: (408) 277-4343\\u000b- Hazardous Waste: (408) 299-7300\\u000b- Recycling Information: <https://www.sanjoseca.gov/index.aspx?NID=243>\",\n  ],\n];\n\nfor (const [find, replace] of replacements) {\n  const results = body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text: \" + find);\n  }\n\n  for (const item of results.items) {\n    item.insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"Level 1 mission complete\" edits to the home utilities contacts\n# document: update the intro heading, refresh provider descriptions /\n# addresses / phone numbers / websites for Sewer, Water and\n# Garbage/Recycling, rename the \"Additional Information\" heading, and swap\n# the old Hazardous Waste / Composting bullets for four new contact lines.\n\n$d = $word.ActiveDocument\n$vtab = [char]11   # <w:br/> renders as a vertical-tab (chr 11) in Range.Text\n\n# Note: we deliberately do NOT use Find.Execute's built-in ReplaceWith\n# parameter here. Doing so routes the inserted text through Word's\n# \"typing\" AutoCorrect pipeline, which silently turns straight apostrophes\n# into curly ones (e.g. \"Jose's\" -> \"Jose's\"). Instead we locate the range\n# with Find.Execute() (no replacement) and assign .Text on the found range\n# directly, which replaces the range contents verbatim.\nfunction Replace-Literal($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = \"\"\n    $found = $find.Execute()\n    if (-not $found) {\n        throw \"Find failed for: $findText\"\n    }\n    $find.Parent.Text = $replaceText\n}\n\nReplace-Literal `\n    \"<b>Utility Providers for San Jose, CA 95148</b>\" `\n    \"<b>Utility Providers in San Jose, CA 95148</b>\"\n\nReplace-Literal `\n    \"- Description: PG&E provides natural gas and electric service to approximately 16 million Californians.\" `\n    \"- Description: PG&E provides natural gas and electric service to approximately 16 million people throughout a 70,000-square-mile service area in Northern and Central California.\"\n\nReplace-Literal `\n    \"- Description: San Jose Water provides water service to over one million people in the greater San Jose metropolitan area.\" `\n    \"- Description: San Jose Water provides essential water services to approximately 1 million people in the greater San Jose metropolitan area.\"\n\nReplace-Literal `\n    \"- Address: 2130 The Alameda, San Jose, CA 95126\" `\n    \"- Address: 110 Paseo De San Antonio, San Jose, CA 95112\"\n\nReplace-Literal `\n    \"- Description: The City of San Jose provides wastewater services to residents and businesses within the city limits.\" `\n    \"- Description: The City of San Jose's Environmental Services department is responsible for wastewater management, recycling, and solid waste services.\"\n\nReplace-Literal `\n    \"- Address: 505 Los Coches St, San Jose, CA 95122\" `\n    \"- Address: 5055 Almaden Expy, San Jose, CA 95118\"\n\n# The Sewer contact number/website line is matched together with the\n# following website text so this occurrence (and not the identical phone\n# number that also appears later, inside the Hazardous Waste bullet) is\n# the one that gets updated.\nReplace-Literal `\n    \"- Contact Number: (408) 535-6000${vtab}- Website: <https://www.sanjoseca.gov/services/wastewater>\" `\n    \"- Contact Number: (408) 277-4343${vtab}- Website: <https://www.sanjoseca.gov/index.aspx?NID=242>\"\n\nReplace-Literal `\n    \"- Description: Republic Services provides residential and commercial waste and recycling services in San Jose.\" `\n    \"- Description: Republic Services provides recycling, trash, and yard waste collection services to residents and businesses in San Jose.\"\n\nReplace-Literal `\n    \"- Address: 1531 Oakland Rd, San Jose, CA 95110\" `\n    \"- Address: 2550 S 10th St, San Jose, CA 95112\"\n\nReplace-Literal `\n    \"- Contact Number: (408) 629-5000\" `\n    \"- Contact Number: (408) 629-8500\"\n\nReplace-Literal `\n    \"<b>Additional Information:</b>\" `\n    \"<b>Additional Contact Information for Garbage/Recycling:</b>\"\n\n# Replace the two old bullet lines (Hazardous Waste + Composting, joined by\n# the line break between them) with the four new contact lines.\nReplace-Literal `\n    \"- <b>Hazardous Waste:</b> The City of San Jose's Household Hazardous Waste Program provides safe disposal of hazardous materials. Contact (408) 535-6000 for more information.${vtab}- <b>Composting:</b> San Jose's GreenCycle Composting Program offers curbside composting services. Visit <https://www.sanjoseca.gov/services/greencycle> for details.\" `\n    \"- Collection Schedule: <https://www.sanjoseca.gov/index.aspx?NID=242>${vtab}- Bulky Item Pickup: (408) 277-4343${vtab}- Hazardous Waste: (408) 299-7300${vtab}- Recycling Information: <https://www.sanjoseca.gov/index.aspx?NID=243>\"\n"}
